$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Make room: insert 4 blank rows starting at row 69. This shifts the old
#    blank rows (69-772) down to (73-776) and leaves the old "end of table"
#    thick-bottom-border row sitting at row 68 (unchanged for now), while the
#    four newly inserted rows become the new (blank) rows 69-72.
# ---------------------------------------------------------------------------
$ws.Range("A69:G72").Insert()

# ---------------------------------------------------------------------------
# 2) Before we touch row 68's formatting, clone its current ("thick bottom
#    border" / last-row-of-table) look onto the brand new row 72, which is
#    going to become the new bottom boundary of the table.
# ---------------------------------------------------------------------------
$ws.Range("B68:E68").Copy()
$ws.Range("B72:E72").PasteSpecial(-4122)
$ws.Range("G68").Copy()
$ws.Range("G72").PasteSpecial(-4122)
$ws.Rows.Item(72).RowHeight = 24.95

# ---------------------------------------------------------------------------
# 3) Clone the regular "in table" look (from row 67) onto row 68 and onto
#    the other 3 newly inserted rows (69-71), so they all look like normal
#    log rows instead of blank/boundary rows.
# ---------------------------------------------------------------------------
$ws.Range("B67:E67").Copy()
$ws.Range("B68:E68").PasteSpecial(-4122)
$ws.Range("B69:E69").PasteSpecial(-4122)
$ws.Range("B70:E70").PasteSpecial(-4122)
$ws.Range("B71:E71").PasteSpecial(-4122)

$ws.Range("G67").Copy()
$ws.Range("G68").PasteSpecial(-4122)
$ws.Range("G69").PasteSpecial(-4122)
$ws.Range("G70").PasteSpecial(-4122)
$ws.Range("G71").PasteSpecial(-4122)

$ws.Rows.Item(68).RowHeight = 24.95
$ws.Rows.Item(69).RowHeight = 24.95
$ws.Rows.Item(70).RowHeight = 24.95
$ws.Rows.Item(71).RowHeight = 24.95

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Fill in the new log entries' data.
#    Row 68: 2020-04-05, 23:19 -> 23:35, "Reviewed report. ..."
#    Row 69: 2020-04-05, 23:35 -> 23:37, "Submitted second attempt to Canvas."
#    Row 70 / 71: new (still) blank entries reserved for future log lines,
#                 just the group number + date filled in.
#    Row 72: new bottom-of-table boundary row, group number + date filled in.
# ---------------------------------------------------------------------------
$ws.Range("B68").Value = 6977
$ws.Range("C68").Value = 43926
$ws.Range("D68").Value = 0.97152777777777777
$ws.Range("E68").Value = 0.98263888888888884
$ws.Range("G68").Value = "Reviewed report. Fixed spelling, formatting and grammer issues. DONE"

$ws.Range("B69").Value = 6977
$ws.Range("C69").Value = 43926
$ws.Range("D69").Value = 0.98263888888888884
$ws.Range("E69").Value = 0.98402777777777783
$ws.Range("G69").Value = "Submitted second attempt to Canvas."

$ws.Range("B70").Value = 6977
$ws.Range("C70").Value = 43926

$ws.Range("B71").Value = 6977
$ws.Range("C71").Value = 43926

$ws.Range("B72").Value = 6977
$ws.Range("C72").Value = 43926

# ---------------------------------------------------------------------------
# 5) Append 4 fresh blank (but correctly-sized) rows at the very end of the
#    sheet, mirroring the long run of blank placeholder rows the log already
#    has.
# ---------------------------------------------------------------------------
$ws.Rows.Item(773).RowHeight = 24.95
$ws.Rows.Item(774).RowHeight = 24.95
$ws.Rows.Item(775).RowHeight = 24.95
$ws.Rows.Item(776).RowHeight = 24.95

# ---------------------------------------------------------------------------
# 6) Move the active selection, matching where the user ended up editing.
# ---------------------------------------------------------------------------
$ws.Range("E71").Select()
